# Scheduled-runner refresh of market-board derived columns (H:N) across the
# per-job Leve tables. Only numeric market-price/profit cells change; leve
# names/levels/EXP/gil/ids (A:G) are untouched. M10 on ALC is cleared because
# that leve no longer has an HQ profit figure.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 6000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 6000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 6000
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = -6586
$ws.Range("H52").Value = 963.1818
$ws.Range("I52").Value = 1182.5
$ws.Range("K52").Value = 3547.5
$ws.Range("M52").Value = -3387.5
$ws.Range("H55").Value = 597
$ws.Range("I55").Value = 194.75
$ws.Range("J55").Value = 999.25
$ws.Range("K55").Value = 194.75
$ws.Range("L55").Value = 999.25
$ws.Range("M55").Value = 19.25
$ws.Range("N55").Value = -1427.25
$ws.Range("H62").Value = 32683.035
$ws.Range("I62").Value = 42356
$ws.Range("K62").Value = 42356
$ws.Range("M62").Value = -41732
$ws.Range("H65").Value = 32683.035
$ws.Range("I65").Value = 42356
$ws.Range("K65").Value = 211780
$ws.Range("M65").Value = -208660
$ws.Range("H100").Value = 2992.32
$ws.Range("I100").Value = 1649
$ws.Range("J100").Value = 5007.3
$ws.Range("K100").Value = 1649
$ws.Range("L100").Value = 5007.3
$ws.Range("M100").Value = -1108
$ws.Range("N100").Value = -6089.3
$ws.Range("H111").Value = 3281.6924
$ws.Range("I111").Value = 3273.2727
$ws.Range("K111").Value = 9819.8181
$ws.Range("M111").Value = -6752.8181
$ws.Range("H132").Value = 1409.4231
$ws.Range("I132").Value = 1353.15
$ws.Range("J132").Value = 1597
$ws.Range("K132").Value = 4059.45
$ws.Range("L132").Value = 4791
$ws.Range("M132").Value = -1529.45
$ws.Range("N132").Value = -9851
$ws.Range("H138").Value = 1490205.8
$ws.Range("I138").Value = 2127.625
$ws.Range("J138").Value = 1667886.8
$ws.Range("K138").Value = 6382.875
$ws.Range("L138").Value = 5003660.4
$ws.Range("M138").Value = -1242.875
$ws.Range("N138").Value = -5013940.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4441.579
$ws.Range("I102").Value = 1199.4166
$ws.Range("J102").Value = 9999.571
$ws.Range("K102").Value = 1199.4166
$ws.Range("L102").Value = 9999.571
$ws.Range("M102").Value = 422.5834
$ws.Range("N102").Value = -13243.571
$ws.Range("H122").Value = 5557464.5
$ws.Range("I122").Value = 7409011.5
$ws.Range("K122").Value = 22227034.5
$ws.Range("M122").Value = -22224584.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2842.7144
$ws.Range("I99").Value = 762.25
$ws.Range("K99").Value = 762.25
$ws.Range("M99").Value = 735.75
$ws.Range("H105").Value = 2261.923
$ws.Range("I105").Value = 1699.5714
$ws.Range("J105").Value = 2918
$ws.Range("K105").Value = 1699.5714
$ws.Range("L105").Value = 2918
$ws.Range("M105").Value = 47.42859999999996
$ws.Range("N105").Value = -6412

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3757.3794
$ws.Range("I31").Value = 1855.5714
$ws.Range("J31").Value = 8749.625
$ws.Range("K31").Value = 1855.5714
$ws.Range("L31").Value = 8749.625
$ws.Range("M31").Value = -1560.5714
$ws.Range("N31").Value = -9339.625
$ws.Range("H34").Value = 3757.3794
$ws.Range("I34").Value = 1855.5714
$ws.Range("J34").Value = 8749.625
$ws.Range("K34").Value = 1855.5714
$ws.Range("L34").Value = 8749.625
$ws.Range("M34").Value = -1653.5714
$ws.Range("N34").Value = -9153.625
$ws.Range("H93").Value = 24188.777
$ws.Range("I93").Value = 23462.375
$ws.Range("K93").Value = 23462.375
$ws.Range("M93").Value = -21590.375
$ws.Range("H141").Value = 455418.66
$ws.Range("J141").Value = 455418.66
$ws.Range("L141").Value = 455418.66
$ws.Range("N141").Value = -465778.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1759953.2
$ws.Range("I4").Value = 760355.5600000001
$ws.Range("J4").Value = 25250500
$ws.Range("K4").Value = 2281066.68
$ws.Range("L4").Value = 75751500
$ws.Range("M4").Value = -2280954.68
$ws.Range("N4").Value = -75751724
$ws.Range("H56").Value = 7710
$ws.Range("I56").Value = 7710
$ws.Range("K56").Value = 7710
$ws.Range("M56").Value = -7180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16676740
$ws.Range("I70").Value = 27787100
$ws.Range("K70").Value = 27787100
$ws.Range("M70").Value = -27786830
$ws.Range("H73").Value = 16676740
$ws.Range("I73").Value = 27787100
$ws.Range("K73").Value = 27787100
$ws.Range("M73").Value = -27786164
$ws.Range("H121").Value = 34500
$ws.Range("J121").Value = 31000
$ws.Range("L121").Value = 31000
$ws.Range("N121").Value = -34494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2513.6155
$ws.Range("I22").Value = 1249.5
$ws.Range("J22").Value = 2743.4546
$ws.Range("K22").Value = 1249.5
$ws.Range("L22").Value = 2743.4546
$ws.Range("M22").Value = -954.5
$ws.Range("N22").Value = -3333.4546
$ws.Range("H27").Value = 2513.6155
$ws.Range("I27").Value = 1249.5
$ws.Range("J27").Value = 2743.4546
$ws.Range("K27").Value = 1249.5
$ws.Range("L27").Value = 2743.4546
$ws.Range("M27").Value = -1142.5
$ws.Range("N27").Value = -2957.4546
$ws.Range("H55").Value = 1363.7727
$ws.Range("I55").Value = 780.5714
$ws.Range("K55").Value = 780.5714
$ws.Range("M55").Value = -607.5714
$ws.Range("H74").Value = 333363400
$ws.Range("J74").Value = 600024600
$ws.Range("L74").Value = 600024600
$ws.Range("N74").Value = -600026596
$ws.Range("H77").Value = 333363400
$ws.Range("J77").Value = 600024600
$ws.Range("L77").Value = 1800073800
$ws.Range("N77").Value = -1800083784

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2330.25
$ws.Range("J81").Value = 3052.75
$ws.Range("L81").Value = 6105.5
$ws.Range("N81").Value = -8227.5
$ws.Range("H84").Value = 2330.25
$ws.Range("J84").Value = 3052.75
$ws.Range("L84").Value = 30527.5
$ws.Range("N84").Value = -41135.5
$ws.Range("H126").Value = 3033.9048
$ws.Range("I126").Value = 2391.5334
$ws.Range("J126").Value = 4639.8335
$ws.Range("K126").Value = 7174.600199999999
$ws.Range("L126").Value = 13919.5005
$ws.Range("M126").Value = -4704.600199999999
$ws.Range("N126").Value = -18859.5005
$ws.Range("H136").Value = 3798.4614
$ws.Range("I136").Value = 2910.6956
$ws.Range("J136").Value = 5074.625
$ws.Range("K136").Value = 8732.086800000001
$ws.Range("L136").Value = 15223.875
$ws.Range("M136").Value = -6182.086800000001
$ws.Range("N136").Value = -20323.875

